# DC-Colos.xlsx — "update generated data" commit
#
# The colo reference table (Sheet1) gained one new row: a Chinese colo
# "XNN" / Xining is inserted right before the existing "IAD" / Ashburn
# row (currently row 271), pushing that row and every row after it
# (through the former last row 330, "YHZ" / Halifax) down by one, so the
# sheet now runs through row 331 and the used-range dimension grows from
# A1:H330 to A1:H331.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 271; this shifts the old row 271 ("IAD", …)
# and everything below it down by one row (old row 330 becomes row 331).
$ws.Rows.Item(271).Insert()

# Populate the newly-opened row 271 with the new colo's data. This table
# has no lat/lon for this colo, so columns G (lat) and H (lon) are left
# blank, matching the other unresolved-coordinate rows elsewhere in the
# sheet (e.g. "TEN"/Tongren, "HYN"/Taizhou just above it).
$a271 = $ws.Cells.Item(271, 1)
$a271.Value = "XNN"
$ws.Cells.Item(271, 2).Value = "Xining, China"
$ws.Cells.Item(271, 3).Value = "Asia"
$ws.Cells.Item(271, 4).Value = "Xining"
$ws.Cells.Item(271, 5).Value = "China"
$ws.Cells.Item(271, 6).Value = "CN"

# Match the header-style formatting used by column A on every other data
# row (bold, thin box border, centered/top-aligned) — the row insert
# above does not automatically carry this over.
$a271.Font.Bold = $true
$a271.Borders.Item(1).LineStyle = 1
$a271.Borders.Item(2).LineStyle = 1
$a271.Borders.Item(3).LineStyle = 1
$a271.Borders.Item(4).LineStyle = 1
$a271.HorizontalAlignment = -4108
$a271.VerticalAlignment = -4160
